# Added LED driver mosfet
# A new MOSFET (SI4896DY-T1-E3) is inserted as the new first part row, the
# previous first part row (MOSFET IRLL2705TRPBF) is replaced by another new
# MOSFET (IRL6372TRPBF), and every other existing part row shifts down by
# one row. A new trailing blank row is appended so the sheet keeps its two
# blank rows at the bottom.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shift all the existing part rows down by one (bottom-up so we never
# clobber a row before it has been read). Row 12 (blank) -> Row 13 (blank),
# row 11 (blank) -> row 12 (blank), row 10 -> 11, ... row 4 -> 5.
for ($r = 12; $r -ge 4; $r--) {
    $dest = $r + 1
    $ws.Range("A$dest").Value2 = $ws.Range("A$r").Value2
    $ws.Range("B$dest").Value2 = $ws.Range("B$r").Value2
    $ws.Range("C$dest").Value2 = $ws.Range("C$r").Value2
    $ws.Range("D$dest").Value2 = $ws.Range("D$r").Value2
    $ws.Range("E$dest").Value2 = $ws.Range("E$r").Value2
    $ws.Range("F$dest").Value2 = $ws.Range("F$r").Value2
    $ws.Range("G$dest").Formula = "=F$dest*E$dest"
}

# --- Row 3 becomes the new MOSFET N-CH 80V 6.7A 8-SOIC part.
$ws.Range("B3").Value2 = "MOSFET N-CH 80V 6.7A 8-SOIC"
$ws.Range("C3").Value2 = "SI4896DY-T1-E3"
$ws.Range("D3").Value2 = "http://www.vishay.com/docs/71300/71300.pdf"
$ws.Range("E3").Value2 = 2.37
$ws.Range("F3").Value2 = 2
$ws.Range("G3").Formula = "=F3*E3"

# --- Row 4 becomes the new MOSFET 2N-CH 30V 8.1A 8SOIC part (replaces the
# old MOSFET N-CH 55V 3.8A SOT223 / IRLL2705TRPBF row that used to live at
# row 3).
$ws.Range("B4").Value2 = "MOSFET 2N-CH 30V 8.1A 8SOIC"
$ws.Range("C4").Value2 = "IRL6372TRPBF"
$ws.Range("D4").Value2 = "http://www.infineon.com/dgdl/irl6372pbf.pdf?fileId=5546d462533600a401535660046e2579"
$ws.Range("E4").Value2 = 1.15
$ws.Range("F4").Value2 = 2
$ws.Range("G4").Formula = "=F4*E4"

# --- New trailing blank row 13 (mirrors the existing blank row 12). Copy
# row 12's formatting (borders/number formats) across first, then clear the
# values and (re)write the G formula, same shape as the other blank row.
$ws.Range("A12:G12").Copy() | Out-Null
$ws.Range("A13:G13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

$ws.Range("A13").Value2 = ""
$ws.Range("B13").Value2 = ""
$ws.Range("C13").Value2 = ""
$ws.Range("D13").Value2 = ""
$ws.Range("E13").Value2 = ""
$ws.Range("F13").Value2 = ""
$ws.Range("G13").Formula = "=F13*E13"

# --- Hyperlinks: this engine's Hyperlinks collection is sheet-wide, so
# rebuild it from scratch once all the cell values are in their final
# places, rather than trying to patch individual entries in place.
$ws.Range("D3").Hyperlinks.Delete() | Out-Null

$ws.Hyperlinks.Add($ws.Range("D3"), "http://www.vishay.com/docs/71300/71300.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D5"), "http://www.st.com/content/ccc/resource/technical/document/datasheet/d4/83/f4/ff/dc/cc/48/7f/CD00001323.pdf/files/CD00001323.pdf/jcr:content/translations/en.CD00001323.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D6"), "http://optoelectronics.liteon.com/upload/download/DS70-2001-006/LTV-355T%20series%20201610.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D7"), "http://www.on-shore.com/wp-content/uploads/2015/09/osttcxx2162.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D8"), "http://optoelectronics.liteon.com/upload/download/DS-22-98-0004/LTST-C150GKT.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D9"), "http://www.molex.com/pdm_docs/sd/039281043_sd.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D10"), "https://media.digikey.com/pdf/Data%20Sheets/Panasonic%20Resistors%20Thermistors%20PDFs/EXBV8V%20Spec.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D11"), "https://media.digikey.com/pdf/Data%20Sheets/Panasonic%20Resistors%20Thermistors%20PDFs/EXBV8V%20Spec.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Range("D4"), "http://www.infineon.com/dgdl/irl6372pbf.pdf?fileId=5546d462533600a401535660046e2579") | Out-Null

# --- Selection / dimension bookkeeping to mirror the saved workbook state.
$ws.Range("A8").Select() | Out-Null
